$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 11.58445689323961
$ws.Range("C2").Value = 8.731790679769469
$ws.Range("E2").Value = 15.99471012893557
$ws.Range("F2").Value = 48.37297769984487
$ws.Range("G2").Value = 3.70816093265667
$ws.Range("I2").Value = 29.92422550733014
$ws.Range("J2").Value = 9.701107968281802
$ws.Range("K2").Value = 12.69504379280688
$ws.Range("B3").Value = 11.41120016731536
$ws.Range("C3").Value = 8.56328953881072
$ws.Range("E3").Value = 15.86921020823669
$ws.Range("F3").Value = 48.040296470916
$ws.Range("G3").Value = 3.711183923550451
$ws.Range("I3").Value = 29.85452460021247
$ws.Range("J3").Value = 9.712920972163992
$ws.Range("K3").Value = 12.58051747831464
$ws.Range("B4").Value = 11.30729525188269
$ws.Range("C4").Value = 8.460985712236932
$ws.Range("E4").Value = 15.79732570899698
$ws.Range("F4").Value = 47.84694589064705
$ws.Range("G4").Value = 3.713134755666194
$ws.Range("I4").Value = 29.81621066816179
$ws.Range("J4").Value = 9.721251046254118
$ws.Range("K4").Value = 12.51339139969444
$ws.Range("B5").Value = 11.2656361037729
$ws.Range("C5").Value = 8.419648591502074
$ws.Range("E5").Value = 15.76936567731931
$ws.Range("F5").Value = 47.77095813298427
$ws.Range("G5").Value = 3.713953637623987
$ws.Range("I5").Value = 29.80172684131909
$ws.Range("J5").Value = 9.724916078445675
$ws.Range("K5").Value = 12.48686951535469
$ws.Range("B6").Value = 11.25876162616991
$ws.Range("C6").Value = 8.412807689028719
$ws.Range("E6").Value = 15.76480437647731
$ws.Range("F6").Value = 47.75851141044101
$ws.Range("G6").Value = 3.714091058559568
$ws.Range("I6").Value = 29.79939007176099
$ws.Range("J6").Value = 9.72554097951843
$ws.Range("K6").Value = 12.48251668824423
$ws.Range("B7").Value = 11.30673057781592
$ws.Range("C7").Value = 8.460426716492165
$ws.Range("E7").Value = 15.79694318985016
$ws.Range("F7").Value = 47.8459096653248
$ws.Range("G7").Value = 3.713145702500267
$ws.Range("I7").Value = 29.81601075919818
$ws.Range("J7").Value = 9.721299379519969
$ws.Range("K7").Value = 12.51303030771139
$ws.Range("B8").Value = 11.52424081902873
$ws.Range("C8").Value = 8.67348839155563
$ws.Range("E8").Value = 15.95038455344651
$ws.Range("F8").Value = 48.25604032582021
$ws.Range("G8").Value = 3.70918365735003
$ws.Range("I8").Value = 29.89926383937571
$ws.Range("J8").Value = 9.704957405539409
$ws.Range("K8").Value = 12.65491050506481
$ws.Range("B9").Value = 11.96771820340764
$ws.Range("C9").Value = 9.09779014578222
$ws.Range("E9").Value = 16.29076268597338
$ws.Range("F9").Value = 49.14410006780663
$ws.Range("G9").Value = 3.702161502955934
$ws.Range("I9").Value = 30.09795898660172
$ws.Range("J9").Value = 9.681469308408712
$ws.Range("K9").Value = 12.95710286530978
$ws.Range("B10").Value = 12.30023557975496
$ws.Range("C10").Value = 9.409957829674656
$ws.Range("E10").Value = 16.56264900599317
$ws.Range("F10").Value = 49.84354764463671
$ws.Range("G10").Value = 3.697452318191553
$ws.Range("I10").Value = 30.26531500407059
$ws.Range("J10").Value = 9.669449275485299
$ws.Range("K10").Value = 13.19180111940954
$ws.Range("B11").Value = 12.45219993935236
$ws.Range("C11").Value = 9.551354391617741
$ws.Range("E11").Value = 16.69055414269653
$ws.Range("F11").Value = 50.17099476288347
$ws.Range("G11").Value = 3.695406493348822
$ws.Range("I11").Value = 30.34602272713071
$ws.Range("J11").Value = 9.665121853919846
$ws.Range("K11").Value = 13.30089293603792
$ws.Range("B12").Value = 12.50978397921362
$ws.Range("C12").Value = 9.604753897203635
$ws.Range("E12").Value = 16.73955026732317
$ws.Range("F12").Value = 50.29623770824637
$ws.Range("G12").Value = 3.694645563601406
$ws.Range("I12").Value = 30.37723520700087
$ws.Range("J12").Value = 9.663647429859941
$ws.Range("K12").Value = 13.34250070110444
$ws.Range("B13").Value = 12.49738148812144
$ws.Range("C13").Value = 9.593260646472645
$ws.Range("E13").Value = 16.72897385590066
$ws.Range("F13").Value = 50.26921039281871
$ws.Range("G13").Value = 3.694808831905656
$ws.Range("I13").Value = 30.37048425803694
$ws.Range("J13").Value = 9.66395766326295
$ws.Range("K13").Value = 13.33352711141948
$ws.Range("B14").Value = 12.45693699276394
$ws.Range("C14").Value = 9.555750807021933
$ws.Range("E14").Value = 16.69457410205575
$ws.Range("F14").Value = 50.18127411100688
$ws.Range("G14").Value = 3.695343615510815
$ws.Range("I14").Value = 30.34857761386995
$ws.Range("J14").Value = 9.664997258387084
$ws.Range("K14").Value = 13.30431029376815
$ws.Range("B15").Value = 12.43216676011639
$ws.Range("C15").Value = 9.532754524299733
$ws.Range("E15").Value = 16.67357496775398
$ws.Range("F15").Value = 50.12757023494878
$ws.Range("G15").Value = 3.695672977965943
$ws.Range("I15").Value = 30.3352435888915
$ws.Range("J15").Value = 9.665655441527029
$ws.Range("K15").Value = 13.28645171949062
$ws.Range("B16").Value = 12.29031308042531
$ws.Range("C16").Value = 9.400700054915836
$ws.Range("E16").Value = 16.55437116531906
$ws.Range("F16").Value = 49.82232730676699
$ws.Range("G16").Value = 3.697587950465808
$ws.Range("I16").Value = 30.26013200437575
$ws.Range("J16").Value = 9.6697550568885
$ws.Range("K16").Value = 13.18471532864474
$ws.Range("B17").Value = 12.20342497024725
$ws.Range("C17").Value = 9.319491520575736
$ws.Range("E17").Value = 16.48229260762448
$ws.Range("F17").Value = 49.63738293872189
$ws.Range("G17").Value = 3.698787356676981
$ws.Range("I17").Value = 30.21522034455111
$ws.Range("J17").Value = 9.672562359048227
$ws.Range("K17").Value = 13.12287262798243
$ws.Range("B18").Value = 12.15351878961106
$ws.Range("C18").Value = 9.272728310200193
$ws.Range("E18").Value = 16.44123521026151
$ws.Range("F18").Value = 49.53188524169911
$ws.Range("G18").Value = 3.699486302652583
$ws.Range("I18").Value = 30.18981965306665
$ws.Range("J18").Value = 9.674284384913246
$ws.Range("K18").Value = 13.0875235825348
$ws.Range("B19").Value = 12.13663528974866
$ws.Range("C19").Value = 9.256887549896328
$ws.Range("E19").Value = 16.42740401573034
$ws.Range("F19").Value = 49.49631882233048
$ws.Range("G19").Value = 3.69972451566795
$ws.Range("I19").Value = 30.18129371801826
$ws.Range("J19").Value = 9.6748858604409
$ws.Range("K19").Value = 13.07559413848924
$ws.Range("B20").Value = 12.21266761574241
$ws.Range("C20").Value = 9.328142336971611
$ws.Range("E20").Value = 16.48992439574187
$ws.Range("F20").Value = 49.65698039223665
$ws.Range("G20").Value = 3.698658738780071
$ws.Range("I20").Value = 30.21995668317189
$ws.Range("J20").Value = 9.672252405659727
$ws.Range("K20").Value = 13.12943326222466
$ws.Range("B21").Value = 12.46881596539654
$ws.Range("C21").Value = 9.566772709164438
$ws.Range("E21").Value = 16.70466328292047
$ws.Range("F21").Value = 50.20707001824958
$ws.Range("G21").Value = 3.695186163288532
$ws.Range("I21").Value = 30.35499454713426
$ws.Range("J21").Value = 9.664687443690314
$ws.Range("K21").Value = 13.31288421888686
$ws.Range("B22").Value = 12.63642026558999
$ws.Range("C22").Value = 9.721866162816802
$ws.Range("E22").Value = 16.84825694333314
$ws.Range("F22").Value = 50.57380532180167
$ws.Range("G22").Value = 3.692996914737528
$ws.Range("I22").Value = 30.44703524510128
$ws.Range("J22").Value = 9.660700937704211
$ws.Range("K22").Value = 13.43449404337552
$ws.Range("B23").Value = 12.54696878661717
$ws.Range("C23").Value = 9.639186982693838
$ws.Range("E23").Value = 16.77133638521775
$ws.Range("F23").Value = 50.37744032999862
$ws.Range("G23").Value = 3.694158039729279
$ws.Range("I23").Value = 30.39756785337361
$ws.Range("J23").Value = 9.662740906354831
$ws.Range("K23").Value = 13.36944418293611
$ws.Range("B24").Value = 12.20848886468819
$ws.Range("C24").Value = 9.324231533220747
$ws.Range("E24").Value = 16.48647287246037
$ws.Range("F24").Value = 49.6481177937566
$ws.Range("G24").Value = 3.698716857671328
$ws.Range("I24").Value = 30.21781407676746
$ws.Range("J24").Value = 9.672392198988264
$ws.Range("K24").Value = 13.12646655761805
$ws.Range("B25").Value = 11.84630858854182
$ws.Range("C25").Value = 8.982680518452542
$ws.Range("E25").Value = 16.19467997166484
$ws.Range("F25").Value = 48.89529168111864
$ws.Range("G25").Value = 3.703981752036057
$ws.Range("I25").Value = 30.04043266212584
$ws.Range("J25").Value = 9.686905138933277
$ws.Range("K25").Value = 12.87298136679454
